$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "297.07"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2.83%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "41.22"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2.40%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.003"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.83%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07521"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "3.09%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.570"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "3.29%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9282"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.03%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1217"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "2.35%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1840"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "6.40%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08898"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "3.05%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04075"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-2.31%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1053"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.07%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001280"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "1.31%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005868"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.75%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.344"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-1.57%"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "1.97%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.10%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.983"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.82%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1420"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "5.12%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.2965"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2.85%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04053"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "4.95%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001265"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.22%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.003878"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.07%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001229"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-4.00%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.02%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02415"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "4.21%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05206"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "4.69%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.005864"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-7.70%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007789"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.45%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1323"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "3.79%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007374"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "0.31%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007835"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "10.85%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.2972"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006238"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-2.88%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.12%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.04514"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "319.75%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.004197"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.10%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002099"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.12%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0001999"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.12%"
